$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to remain Text (matches the
# source data, which stores prices/percentages as inline strings, even
# when they look numeric), without leaving a stray cell-style behind.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "30.614.46"
Set-TextValue "E2" "  +0.49%  "
Set-TextValue "D3" "1.881.90"
Set-TextValue "E3" "  +0.19%  "
Set-TextValue "D4" "1.001"
Set-TextValue "D5" "249.56"
Set-TextValue "E5" "  +1.04%  "
Set-TextValue "E6" "  +0.04%  "
Set-TextValue "D7" "0.4756"
Set-TextValue "E7" "  -0.10%  "
Set-TextValue "D8" "0.2935"
Set-TextValue "E8" "  +1.17%  "
Set-TextValue "D9" "0.06529"
Set-TextValue "E9" "  +0.10%  "
Set-TextValue "D10" "21.93"
Set-TextValue "E10" "  +1.60%  "
Set-TextValue "D11" "0.07737"
Set-TextValue "E11" "  +0.02%  "
Set-TextValue "D12" "96.96"
Set-TextValue "E12" "  -0.05%  "
Set-TextValue "D13" "0.7393"
Set-TextValue "E13" "  -0.81%  "
Set-TextValue "D14" "1.882.33"
Set-TextValue "E14" "  +0.16%  "
Set-TextValue "D15" "5.257"
Set-TextValue "E15" "  +2.67%  "
Set-TextValue "D16" "274.47"
Set-TextValue "E16" "  +0.16%  "
Set-TextValue "D17" "30.604.40"
Set-TextValue "E17" "  +0.51%  "
Set-TextValue "D18" "13.17"
Set-TextValue "E18" "  -3.35%  "
Set-TextValue "D19" "0.000007538"
Set-TextValue "E19" "  -0.41%  "
Set-TextValue "D21" "2.129.91"
Set-TextValue "E21" "  +0.18%  "
Set-TextValue "D22" "5.316"
Set-TextValue "E22" "  +1.29%  "
Set-TextValue "E23" "  +0.05%  "
Set-TextValue "D24" "6.235"
Set-TextValue "E24" "  +1.06%  "
Set-TextValue "D25" "9.213"
Set-TextValue "E25" "  -0.82%  "
Set-TextValue "D26" "163.94"
Set-TextValue "E26" "  -0.23%  "
Set-TextValue "D27" "18.86"
Set-TextValue "E27" "  -0.10%  "
Set-TextValue "D28" "1.916"
Set-TextValue "E28" "  -2.17%  "
Set-TextValue "E29" "  -2.06%  "
Set-TextValue "D30" "0.09689"
Set-TextValue "E30" "  -3.07%  "
Set-TextValue "E31" "  -0.40%  "
Set-TextValue "D32" "4.297"
Set-TextValue "E32" "  -0.68%  "
Set-TextValue "D33" "4.146"
Set-TextValue "E33" "  +2.04%  "
Set-TextValue "D34" "0.04865"
Set-TextValue "E34" "  +1.88%  "
Set-TextValue "D35" "1.127"
Set-TextValue "E35" "  +0.27%  "
Set-TextValue "D36" "0.6987"
Set-TextValue "E36" "  -0.09%  "
Set-TextValue "D37" "2.719"
Set-TextValue "E37" "  +0.09%  "
Set-TextValue "D38" "0.01902"
Set-TextValue "E38" "  +1.85%  "
Set-TextValue "D39" "2.774"
Set-TextValue "E39" "  +1.43%  "
Set-TextValue "D40" "6.319"
Set-TextValue "E40" "  -0.63%  "
Set-TextValue "D41" "74.78"
Set-TextValue "E41" "  +6.68%  "
Set-TextValue "D42" "2.015"
Set-TextValue "E42" "  +4.27%  "
Set-TextValue "D43" "0.4240"
Set-TextValue "E43" "  +1.72%  "
Set-TextValue "B44" "TrustWalletToken"
Set-TextValue "C44" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D44" "0.8407"
Set-TextValue "E44" "  +0.23%  "
Set-TextValue "B45" "PaxDollar"
Set-TextValue "C45" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D45" "1.000"
Set-TextValue "E45" "  +0.05%  "
Set-TextValue "D46" "102.69"
Set-TextValue "E46" "  +0.02%  "
Set-TextValue "D47" "9.410"
Set-TextValue "E47" "  +1.12%  "
Set-TextValue "D48" "7.056"
Set-TextValue "E48" "  -0.32%  "
Set-TextValue "D49" "35.64"
Set-TextValue "E49" "  +0.84%  "
Set-TextValue "D50" "916.55"
Set-TextValue "E50" "  -1.03%  "
Set-TextValue "E51" "  +2.21%  "
